$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.145.94'
$ws.Range("E2").Value = '  -0.48%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.656.99'
$ws.Range("E3").Value = '  -0.44%  '
$ws.Range("E4").Value = '  -0.52%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '218.77'
$ws.Range("E5").Value = '  +0.23%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5248'
$ws.Range("E6").Value = '  -1.48%  '
$ws.Range("E7").Value = '  -0.50%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2621'
$ws.Range("E8").Value = '  -0.63%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06294'
$ws.Range("E9").Value = '  -1.06%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.59'
$ws.Range("E10").Value = '  +0.32%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07803'
$ws.Range("E11").Value = '  -0.26%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.497'
$ws.Range("E12").Value = '  -1.50%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.658.98'
$ws.Range("E13").Value = '  -0.43%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.884.03'
$ws.Range("E14").Value = '  -0.43%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5557'
$ws.Range("E15").Value = '  +0.44%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0₅8008'
$ws.Range("E16").Value = '  -2.46%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '65.05'
$ws.Range("E17").Value = '  -0.96%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '26.154.45'
$ws.Range("E18").Value = '  -0.53%  '
$ws.Range("E19").Value = '  -0.49%  '
$ws.Range("E20").Value = '  -0.97%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '195.23'
$ws.Range("E21").Value = '  +0.97%  '
$ws.Range("E22").Value = '  -0.87%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.967'
$ws.Range("E23").Value = '  -1.18%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '146.68'
$ws.Range("E25").Value = '  +0.47%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1206'
$ws.Range("E26").Value = '  -1.76%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.168'
$ws.Range("E27").Value = '  -0.34%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.96'
$ws.Range("E28").Value = '  -0.72%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.501'
$ws.Range("E29").Value = '  +0.99%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05737'
$ws.Range("E30").Value = '  -2.25%  '
$ws.Range("E31").Value = '  -0.69%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.485'
$ws.Range("E32").Value = '  -2.93%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.345'
$ws.Range("E33").Value = '  +2.10%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.589'
$ws.Range("E34").Value = '  -1.19%  '
$ws.Range("E35").Value = '  -0.63%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9532'
$ws.Range("E36").Value = '  -0.81%  '
$ws.Range("E37").Value = '  -0.24%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5708'
$ws.Range("E38").Value = '  -1.56%  '
$ws.Range("E39").Value = '  -0.43%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.950'
$ws.Range("E40").Value = '  +1.90%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.062.18'
$ws.Range("E41").Value = '  +1.35%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8448'
$ws.Range("E42").Value = '  -2.19%  '
$ws.Range("E43").Value = '  -0.45%  '
$ws.Range("E44").Value = '  -0.64%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.794.36'
$ws.Range("E45").Value = '  -0.45%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '57.89'
$ws.Range("E46").Value = '  +0.49%  '
$ws.Range("E47").Value = '  +2.48%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.05390'
$ws.Range("E48").Value = '  +4.44%  '
$ws.Range("E49").Value = '  -0.42%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4399'
$ws.Range("E50").Value = '  +0.44%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.009'
$ws.Range("E51").Value = '  -0.37%  '
